# Apply updated crypto price/volume/coin data to Sheet1 (cells B2:E51).
# Generated from the target OOXML diff: each entry is a cell reference and its
# new text value. Values that look like plain decimal numbers are written with
# NumberFormat '@' (Text) first, then the style is reset to 'Normal' so the
# underlying value stays a text string (matching the source inlineStr cells)
# without leaving a lingering custom number-format style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.605.19' }
    @{ Cell = 'E2'; Value = '  +0.43%  ' }
    @{ Cell = 'D3'; Value = '1.738.97' }
    @{ Cell = 'E3'; Value = '  +0.61%  ' }
    @{ Cell = 'E4'; Value = '  +0.12%  ' }
    @{ Cell = 'D5'; Value = '245.95' }
    @{ Cell = 'E5'; Value = '  +0.22%  ' }
    @{ Cell = 'E6'; Value = '  +0.08%  ' }
    @{ Cell = 'D7'; Value = '0.4923' }
    @{ Cell = 'E7'; Value = '  +2.66%  ' }
    @{ Cell = 'D8'; Value = '0.2672' }
    @{ Cell = 'E8'; Value = '  -0.56%  ' }
    @{ Cell = 'D9'; Value = '0.06274' }
    @{ Cell = 'E9'; Value = '  +0.77%  ' }
    @{ Cell = 'D10'; Value = '1.752.40' }
    @{ Cell = 'E10'; Value = '  +1.37%  ' }
    @{ Cell = 'D11'; Value = '0.07046' }
    @{ Cell = 'E11'; Value = '  -1.20%  ' }
    @{ Cell = 'D12'; Value = '15.75' }
    @{ Cell = 'E12'; Value = '  +0.04%  ' }
    @{ Cell = 'D13'; Value = '0.6139' }
    @{ Cell = 'E13'; Value = '  -0.89%  ' }
    @{ Cell = 'D14'; Value = '4.580' }
    @{ Cell = 'E14'; Value = '  +1.08%  ' }
    @{ Cell = 'D16'; Value = '1.000' }
    @{ Cell = 'E16'; Value = '  +0.07%  ' }
    @{ Cell = 'D17'; Value = '26.620.38' }
    @{ Cell = 'E17'; Value = '  +0.47%  ' }
    @{ Cell = 'D18'; Value = '1.000' }
    @{ Cell = 'E18'; Value = '  +0.08%  ' }
    @{ Cell = 'D19'; Value = '0.000007261' }
    @{ Cell = 'E19'; Value = '  +4.33%  ' }
    @{ Cell = 'E20'; Value = '  -1.33%  ' }
    @{ Cell = 'D21'; Value = '1.969.98' }
    @{ Cell = 'E21'; Value = '  +0.89%  ' }
    @{ Cell = 'D22'; Value = '4.563' }
    @{ Cell = 'E22'; Value = '  +0.49%  ' }
    @{ Cell = 'D23'; Value = '8.706' }
    @{ Cell = 'E23'; Value = '  -2.66%  ' }
    @{ Cell = 'D24'; Value = '5.272' }
    @{ Cell = 'E24'; Value = '  -0.53%  ' }
    @{ Cell = 'D25'; Value = '139.14' }
    @{ Cell = 'E25'; Value = '  +1.98%  ' }
    @{ Cell = 'D26'; Value = '15.41' }
    @{ Cell = 'E26'; Value = '  +0.31%  ' }
    @{ Cell = 'D27'; Value = '1.421' }
    @{ Cell = 'E27'; Value = '  +1.11%  ' }
    @{ Cell = 'D28'; Value = '1.758' }
    @{ Cell = 'E28'; Value = '  -2.53%  ' }
    @{ Cell = 'D29'; Value = '107.44' }
    @{ Cell = 'E29'; Value = '  +0.59%  ' }
    @{ Cell = 'D30'; Value = '4.024' }
    @{ Cell = 'E30'; Value = '  +0.99%  ' }
    @{ Cell = 'E31'; Value = '  +0.31%  ' }
    @{ Cell = 'D32'; Value = '3.726' }
    @{ Cell = 'E32'; Value = '  -0.18%  ' }
    @{ Cell = 'D33'; Value = '0.04619' }
    @{ Cell = 'E33'; Value = '  +1.14%  ' }
    @{ Cell = 'B34'; Value = 'HuobiToken' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Cell = 'D34'; Value = '2.612' }
    @{ Cell = 'E34'; Value = '  -0.20%  ' }
    @{ Cell = 'B35'; Value = 'ARBITRUM' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'D35'; Value = '1.014' }
    @{ Cell = 'E35'; Value = '  +2.42%  ' }
    @{ Cell = 'B36'; Value = 'ImmutableX' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D36'; Value = '0.6387' }
    @{ Cell = 'E36'; Value = '  +0.12%  ' }
    @{ Cell = 'B37'; Value = 'RenderToken' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D37'; Value = '2.065' }
    @{ Cell = 'E37'; Value = '  -0.98%  ' }
    @{ Cell = 'D38'; Value = '0.9055' }
    @{ Cell = 'E38'; Value = '  -3.11%  ' }
    @{ Cell = 'B39'; Value = 'MXToken' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D39'; Value = '2.426' }
    @{ Cell = 'E39'; Value = '  +0.44%  ' }
    @{ Cell = 'B40'; Value = 'PaxDollar' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' }
    @{ Cell = 'D40'; Value = '1.003' }
    @{ Cell = 'E40'; Value = '  -0.32%  ' }
    @{ Cell = 'B41'; Value = 'VeChain' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D41'; Value = '0.01505' }
    @{ Cell = 'E41'; Value = '  +0.31%  ' }
    @{ Cell = 'B42'; Value = 'Quant' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' }
    @{ Cell = 'D42'; Value = '101.94' }
    @{ Cell = 'E42'; Value = '  -3.87%  ' }
    @{ Cell = 'B43'; Value = 'FraxShare' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D43'; Value = '5.428' }
    @{ Cell = 'E43'; Value = '  -5.04%  ' }
    @{ Cell = 'B44'; Value = 'TheSandbox' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Cell = 'D44'; Value = '0.3932' }
    @{ Cell = 'E44'; Value = '  +0.39%  ' }
    @{ Cell = 'B45'; Value = 'Aptos' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'D45'; Value = '6.874' }
    @{ Cell = 'E45'; Value = '  -1.61%  ' }
    @{ Cell = 'B46'; Value = 'Algorand' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Cell = 'D46'; Value = '0.1186' }
    @{ Cell = 'E46'; Value = '  -0.45%  ' }
    @{ Cell = 'B47'; Value = 'Cronos' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'D47'; Value = '0.05390' }
    @{ Cell = 'E47'; Value = '  +1.33%  ' }
    @{ Cell = 'B48'; Value = 'Elrond' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' }
    @{ Cell = 'D48'; Value = '30.61' }
    @{ Cell = 'E48'; Value = '  -1.51%  ' }
    @{ Cell = 'B49'; Value = 'EnergySwap' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D49'; Value = '7.800' }
    @{ Cell = 'E49'; Value = '  -1.56%  ' }
    @{ Cell = 'B50'; Value = 'NEARProtocol' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = 'D50'; Value = '1.255' }
    @{ Cell = 'E50'; Value = '  -1.08%  ' }
    @{ Cell = 'B51'; Value = 'Aave' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D51'; Value = '51.80' }
    @{ Cell = 'E51'; Value = '  +0.75%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Value -match '^-?\d+(\.\d+)?$') {
        # Plain numeric-looking text (e.g. '245.95') -- force Text format so
        # Excel stores it as a string instead of auto-converting to a number,
        # then reset the style so no stray number-format style is left behind.
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $u.Value
    }
}
